$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 16.73645575277018
    "D2" = 11.73410725120842
    "E2" = 18.10012553535841
    "F2" = 29.14199549063542
    "G2" = 27.48748284503913
    "H2" = 13.78297157516984
    "J2" = 12.33477196596838
    "L2" = 8.920168514362683
    "M2" = 15.39100559533579
    "N2" = 19.45251235619546
    "O2" = 20.84525897907064
    "B3" = 16.53283508584281
    "D3" = 11.7667170440744
    "E3" = 18.12424347219878
    "F3" = 29.19100092302396
    "G3" = 27.35491365961642
    "H3" = 13.80982709129547
    "J3" = 12.33918486265699
    "L3" = 8.783460466770141
    "M3" = 15.26942218383456
    "N3" = 19.4788699560045
    "O3" = 20.85713276023275
    "B4" = 16.40878950349778
    "D4" = 11.78788477026587
    "E4" = 18.14053426271277
    "F4" = 29.22826030217788
    "G4" = 27.28250594985066
    "H4" = 13.82876276432676
    "J4" = 12.34257678219054
    "L4" = 8.698890822257509
    "M4" = 15.19573019687533
    "N4" = 19.49717136066739
    "O4" = 20.86938129633954
    "B5" = 16.35853878715943
    "D5" = 11.79679950932524
    "E5" = 18.14754677792251
    "F5" = 29.24524355731175
    "G5" = 27.25528690026858
    "H5" = 13.83709397518701
    "J5" = 12.34413140226814
    "L5" = 8.664304754567866
    "M5" = 15.16596519889268
    "N5" = 19.50516289960403
    "O5" = 20.87561897435698
    "B6" = 16.35021428199363
    "D6" = 11.7982972550542
    "E6" = 18.14873381620735
    "F6" = 29.2481722103026
    "G6" = 27.25090608056239
    "H6" = 13.83851448456749
    "J6" = 12.34439998048774
    "L6" = 8.658555361903158
    "M6" = 15.16103943216011
    "N6" = 19.50652214492904
    "O6" = 20.87672998502585
    "B7" = 16.4081105251263
    "D7" = 11.78800382753092
    "E7" = 18.14062732056379
    "F7" = 29.22848206179384
    "G7" = 27.28212956873957
    "H7" = 13.82887263332551
    "J7" = 12.34259704931189
    "L7" = 8.69842483449259
    "M7" = 15.19532767167691
    "N7" = 19.49727697559176
    "O7" = 20.86946037465553
    "B8" = 16.66607609236047
    "D8" = 11.74511387953175
    "E8" = 18.10813444780944
    "F8" = 29.15740274119919
    "G8" = 27.43992442673188
    "H8" = 13.79172343724286
    "J8" = 12.3361523821658
    "L8" = 8.873183364502522
    "M8" = 15.34889846333627
    "N8" = 19.46116155551773
    "O8" = 20.84832395047616
    "B9" = 17.17737912410628
    "D9" = 11.67005901000644
    "E9" = 18.05612626594493
    "F9" = 29.07502257346488
    "G9" = 27.81932331969288
    "H9" = 13.73829781908943
    "J9" = 12.3288941337992
    "L9" = 9.209369997413393
    "M9" = 15.65662441580921
    "N9" = 19.4070965435142
    "O9" = 20.84620954424974
    "B10" = 17.55323898444518
    "D10" = 11.62038673089277
    "E10" = 18.02498579318376
    "F10" = 29.04937287849345
    "G10" = 28.13863183030236
    "H10" = 13.71090205959559
    "J10" = 12.32679484454222
    "L10" = 9.450432509690522
    "M10" = 15.8852969054409
    "N10" = 19.3775286342862
    "O10" = 20.86859199346652
    "B11" = 17.72362394845359
    "D11" = 11.59896723638887
    "E11" = 18.01233955052062
    "F11" = 29.04528898487953
    "G11" = 28.29218730582804
    "H11" = 13.70101462175482
    "J11" = 12.32653191567819
    "L11" = 9.55841543031808
    "M11" = 15.98958219769665
    "N11" = 19.36626852275334
    "O11" = 20.88394895330197
    "B12" = 17.7880057623545
    "D12" = 11.59102466956632
    "E12" = 18.00776805093908
    "F12" = 29.04483282004717
    "G12" = 28.35148078102802
    "H12" = 13.69764068484265
    "J12" = 12.32653100051839
    "L12" = 9.599033275996211
    "M12" = 16.02908413674659
    "N12" = 19.36231839718539
    "O12" = 20.89050527598042
    "B13" = 17.77414702027652
    "D13" = 11.59272775864419
    "E13" = 18.00874295692646
    "F13" = 29.04488258400023
    "G13" = 28.33866072529333
    "H13" = 13.69835085911404
    "J13" = 12.32652682401448
    "L13" = 9.590298124867722
    "M13" = 16.02057659711108
    "N13" = 19.36315518827984
    "O13" = 20.88906035779853
    "B14" = 17.72892377259554
    "D14" = 11.59831042286045
    "E14" = 18.01195910108383
    "F14" = 29.04522960875358
    "G14" = 28.29704274647758
    "H14" = 13.70072962693766
    "J14" = 12.32652986828456
    "L14" = 9.561762731590523
    "M14" = 15.99283196764656
    "N14" = 19.36593726057033
    "O14" = 20.88447351644482
    "B15" = 17.70120352160957
    "D15" = 11.601751893709
    "E15" = 18.01395735364271
    "F15" = 29.04558414110818
    "G15" = 28.27169821609819
    "H15" = 13.70223490072384
    "J15" = 12.32654455429612
    "L15" = 9.544247517629955
    "M15" = 15.97583828647094
    "N15" = 19.36768219583564
    "O15" = 20.88176033380772
    "B16" = 17.54208732684557
    "D16" = 11.62181016521264
    "E16" = 18.02584273519824
    "F16" = 29.04979237806574
    "G16" = 28.12875957965688
    "H16" = 13.71160004020328
    "J16" = 12.32682587961823
    "L16" = 9.44333907129659
    "M16" = 15.87848462839128
    "N16" = 19.37830848923437
    "O16" = 20.86769225866071
    "B17" = 17.44428358062749
    "D17" = 11.63441615837895
    "E17" = 18.03352243413419
    "F17" = 29.05431654521034
    "G17" = 28.04316380130034
    "H17" = 13.7180047675679
    "J17" = 12.32717510152936
    "L17" = 9.380982323246018
    "M17" = 15.81880950600713
    "N17" = 19.38538760005061
    "O17" = 20.86038535999232
    "B18" = 17.3879765261158
    "D18" = 11.64177758627717
    "E18" = 18.03808270745438
    "F18" = 29.05763266062959
    "G18" = 27.99471571176109
    "H18" = 13.72193098256267
    "J18" = 12.32744115605813
    "L18" = 9.344960273695053
    "M18" = 15.78451208379518
    "N18" = 19.38966556671882
    "O18" = 20.85666996283564
    "B19" = 17.3689044905773
    "D19" = 11.64428908887996
    "E19" = 18.03965135318692
    "F19" = 29.05887805564678
    "G19" = 27.97844815935899
    "H19" = 13.72330196148457
    "J19" = 12.32754246040949
    "L19" = 9.33273800597776
    "M19" = 15.77290485343145
    "N19" = 19.39114947085489
    "O19" = 20.85549577978502
    "B20" = 17.4547008268169
    "D20" = 11.63306276702524
    "E20" = 18.03269011383845
    "F20" = 29.05376105007761
    "G20" = 28.05219478327232
    "H20" = 13.71729788754078
    "J20" = 12.32713118543933
    "L20" = 9.387636695175956
    "M20" = 15.82515951035149
    "N20" = 19.38461267883322
    "O20" = 20.86111277510919
    "B21" = 17.74221113711787
    "D21" = 11.59666609051051
    "E21" = 18.0110085519617
    "F21" = 29.04509809442106
    "G21" = 28.30923626817389
    "H21" = 13.70002087856494
    "J21" = 12.32652630371959
    "L21" = 9.570151921423891
    "M21" = 16.00098113730973
    "N21" = 19.3651115904944
    "O21" = 20.8858007019986
    "B22" = 17.9292800870643
    "D22" = 11.57386077742357
    "E22" = 17.99810477496361
    "F22" = 29.04579087427189
    "G22" = 28.48387962923858
    "H22" = 13.69088717386358
    "J22" = 12.32670557864932
    "L22" = 9.687831575115844
    "M22" = 16.11594495091802
    "N22" = 19.35419520853594
    "O22" = 20.90625268575016
    "B23" = 17.82953168329494
    "D23" = 11.58594276822543
    "E23" = 18.00487628071415
    "F23" = 29.04483998108985
    "G23" = 28.39007710329715
    "H23" = 13.69556461512465
    "J23" = 12.32655760906517
    "L23" = 9.625180649166468
    "M23" = 16.0545902412576
    "N23" = 19.35985455041504
    "O23" = 20.89494330731488
    "B24" = 17.44999142982352
    "D24" = 11.63367427988735
    "E24" = 18.03306595375646
    "F24" = 29.05400996189475
    "G24" = 28.0481094989131
    "H24" = 13.71761670778362
    "J24" = 12.32715083649364
    "L24" = 9.384628787951495
    "M24" = 15.82228863835734
    "N24" = 19.38496237244095
    "O24" = 20.86078239839372
    "B25" = 17.03879897700814
    "D25" = 11.68939921947856
    "E25" = 18.0689492843708
    "F25" = 29.09119029266744
    "G25" = 27.70939766966956
    "H25" = 13.75066947724375
    "J25" = 12.33028648435531
    "L25" = 9.209369997413393
    "M25" = 15.65662441580921
    "N25" = 19.40709654351469
    "O25" = 20.84257237958606
}

foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

Write-Host "Updated $($values.Count) cells"